$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Contest 44 - SRH vs CSK (row 56) - fill in the 9 players' raw scores
$ws.Range("E56").Value = 60
$ws.Range("H56").Value = 80
$ws.Range("K56").Value = 50
$ws.Range("N56").Value = 30
$ws.Range("Q56").Value = 40
$ws.Range("T56").Value = 20
$ws.Range("W56").Value = 0
$ws.Range("Z56").Value = 70
$ws.Range("AC56").Value = 100
